$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "2024-11-21 22:11:45"
$ws.Range("B8").Value = 20
$ws.Range("C8").Value = "Alerta Amarelo, Chuvas Intensas"

$ws.Range("A9").Value = "2024-11-22 12:22:29"
$ws.Range("B9").Value = 21

$ws.Range("A10").Value = "2024-11-22 12:23:33"
$ws.Range("B10").Value = 21
$ws.Range("C10").Value = "Alerta Umidade nao encontrado"
